# Update defuzzification weights (column D, "Skor Kelayakan") with new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = 100
$ws.Range("D3").Value  = 98.75
$ws.Range("D4").Value  = 95
$ws.Range("D5").Value  = 93.75
$ws.Range("D6").Value  = 92.5
$ws.Range("D7").Value  = 88.7175
$ws.Range("D8").Value  = 87.5
$ws.Range("D9").Value  = 86.74000000000001
$ws.Range("D10").Value = 83.75
$ws.Range("D11").Value = 78.75
